# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45188 (2023-09-19) to 45189 (2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Set the new date serial value across the whole C2:C<lastRow> range at once.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45189
